$d = $word.ActiveDocument

$replacements = @(
    @("647÷8=80, 7", "722÷7=103, 1"),
    @("624÷9=69, 3", "730÷5=146, 0"),
    @("364÷2=182, 0", "503÷6=83, 5"),
    @("198÷5=39, 3", "655÷3=218, 1"),
    @("383÷5=76, 3", "556÷4=139, 0"),
    @("918÷7=131, 1", "660÷9=73, 3"),
    @("557÷8=69, 5", "770÷8=96, 2"),
    @("285÷5=57, 0", "293÷7=41, 6"),
    @("714÷3=238, 0", "330÷7=47, 1"),
    @("571÷2=285, 1", "142÷2=71, 0"),
    @("947÷6=157, 5", "618÷4=154, 2"),
    @("603÷6=100, 3", "726÷7=103, 5"),
    @("203÷7=29, 0", "749÷6=124, 5"),
    @("971÷9=107, 8", "878÷4=219, 2"),
    @("806÷9=89, 5", "947÷8=118, 3"),
    @("216÷7=30, 6", "201÷4=50, 1"),
    @("513÷6=85, 3", "568÷9=63, 1"),
    @("651÷5=130, 1", "698÷2=349, 0"),
    @("389÷3=129, 2", "996÷8=124, 4"),
    @("693÷8=86, 5", "120÷3=40, 0"),
    @("484÷3=161, 1", "144÷2=72, 0"),
    @("495÷3=165, 0", "430÷6=71, 4"),
    @("203÷4=50, 3", "890÷6=148, 2"),
    @("697÷5=139, 2", "814÷7=116, 2"),
    @("961÷2=480, 1", "123÷6=20, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
